# Camera Config file added
#
# Updates the "Test2_2" sheet (the 3rd / active sheet): flips the
# video_file / gz_pose_file flags on, swaps the second marker's SDF file
# for the new "DICT_4X4_50_s1000_id6.sdf" config, and leaves a new blank,
# formatted cell a few rows below the table (A14) where editing stopped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test2_2")

# video_file: now enabled
$ws.Range("B3").Value = 1

# gz_pose_file: now enabled
$ws.Range("B4").Value = 1

# second marker entry now points at the new camera config file
$ws.Range("B8").Value = "DICT_4X4_50_s1000_id6.sdf"

# touch formatting on the new cell below the table so it carries its own
# (blank) style, mirroring where the cursor ended up after editing
$ws.Range("A14").WrapText = $true

# cursor ends up on D1 after the edit
$ws.Range("D1").Select()
